# Updated symbol list with refreshed Price (col D) and Volume(1h) (col E)
# percentages, mirroring the upstream GitHub Actions data-refresh commit.
#
# Every D/E cell in this sheet is stored as literal text (e.g. "299.71",
# "-0.51%") rather than as a number/percentage, so a plain .Value = "..."
# assignment would let Excel auto-coerce the text into a real number or
# percent (and stamp a numeric/percent style on the cell). To keep the
# cell a genuine text value with its original (default) style, each cell is
# briefly switched to the Text number format ("@") before the write and
# restored to its original style immediately after.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = $originalStyle
}

Set-TextValue "D2" "299.71"
Set-TextValue "E2" "-0.51%"
Set-TextValue "D3" "31.72"
Set-TextValue "E3" "1.01%"
Set-TextValue "D4" "5.144"
Set-TextValue "E4" "0.81%"
Set-TextValue "D5" "0.08101"
Set-TextValue "E5" "9.99%"
Set-TextValue "D6" "2.565"
Set-TextValue "E6" "17.38%"
Set-TextValue "E7" "-1.67%"
Set-TextValue "D8" "3.909"
Set-TextValue "E8" "2.16%"
Set-TextValue "D9" "0.9322"
Set-TextValue "E9" "1.61%"
Set-TextValue "E10" "3.41%"
Set-TextValue "D11" "0.07345"
Set-TextValue "E11" "-1.82%"
Set-TextValue "D12" "0.08883"
Set-TextValue "E12" "9.37%"
Set-TextValue "D13" "0.03029"
Set-TextValue "E13" "0.25%"
Set-TextValue "D14" "0.09993"
Set-TextValue "E14" "0.78%"
Set-TextValue "D15" "0.001510"
Set-TextValue "E15" "0.92%"
Set-TextValue "D16" "0.005806"
Set-TextValue "E16" "-5.69%"
Set-TextValue "D17" "3.566"
Set-TextValue "E17" "3.12%"
Set-TextValue "D18" "2.284"
Set-TextValue "E18" "2.67%"
Set-TextValue "E19" "-0.25%"
Set-TextValue "E20" "-0.60%"
Set-TextValue "D21" "4.186"
Set-TextValue "E21" "-10.11%"
Set-TextValue "E22" "7.23%"
Set-TextValue "D23" "0.04633"
Set-TextValue "E23" "-0.05%"
Set-TextValue "D24" "0.001240"
Set-TextValue "E24" "1.15%"
Set-TextValue "D25" "0.004528"
Set-TextValue "E25" "1.18%"
Set-TextValue "D26" "0.0001201"
Set-TextValue "D27" "0.0003409"
Set-TextValue "E27" "-0.52%"
Set-TextValue "D39" "0.01761"
Set-TextValue "E39" "1.89%"
Set-TextValue "E40" "2.37%"
Set-TextValue "D41" "0.006925"
Set-TextValue "E41" "-4.49%"
Set-TextValue "D42" "0.1375"
Set-TextValue "E42" "1.98%"
Set-TextValue "D43" "0.002212"
Set-TextValue "E43" "-0.83%"
Set-TextValue "D44" "0.01035"
Set-TextValue "E44" "-2.71%"
Set-TextValue "D45" "0.00006323"
Set-TextValue "E45" "0.31%"
Set-TextValue "E46" "-0.04%"
Set-TextValue "D47" "0.008397"
Set-TextValue "E47" "-15.96%"
Set-TextValue "D48" "0.7486"
Set-TextValue "E48" "-8.84%"
Set-TextValue "D49" "0.00002099"
Set-TextValue "E49" "-0.04%"
Set-TextValue "D50" "0.0001999"
Set-TextValue "E50" "0.03%"
